$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.494.31"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "3.007.79"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'583.81"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'146.11"
$ws.Range("E6").Value = "  -5.27%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").Value = "3.007.10"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("E10").Value = "  -4.30%  "

$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("E13").Value = "  -3.42%  "

$ws.Range("D14").Value = "'34.63"
$ws.Range("E14").Value = "  -6.04%  "

$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("D16").Value = "3.502.81"
$ws.Range("E16").Value = "  -1.91%  "

$ws.Range("D17").Value = "'7.07"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").Value = "62.444.62"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").Value = "3.010.29"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'459.03"
$ws.Range("E20").Value = "  -5.42%  "

$ws.Range("D21").Value = "'13.95"
$ws.Range("E21").Value = "  -3.02%  "

$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").Value = "'7.41"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("D24").Value = "'81.61"
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'12.35"
$ws.Range("E25").Value = "  -4.38%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  -9.06%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  -6.20%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'2.61"
$ws.Range("E30").Value = "  -2.94%  "

$ws.Range("D31").Value = "'7.02"
$ws.Range("E31").Value = "  -4.76%  "

$ws.Range("E32").Value = "  -6.22%  "

$ws.Range("D33").Value = "'28.06"
$ws.Range("E33").Value = "  +1.99%  "

$ws.Range("E34").Value = "  -2.77%  "

$ws.Range("D35").Value = "0.0₃0808"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  -3.80%  "

$ws.Range("D37").Value = "'5.76"
$ws.Range("E37").Value = "  -3.86%  "

$ws.Range("D38").Value = "'2.12"
$ws.Range("E38").Value = "  -5.74%  "

$ws.Range("D39").Value = "'50.37"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").Value = "'9.13"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  -13.65%  "

$ws.Range("D42").Value = "'0.120"
$ws.Range("E42").Value = "  +4.56%  "

$ws.Range("D43").Value = "'390.09"
$ws.Range("E43").Value = "  -11.43%  "

$ws.Range("D44").Value = "'0.0358"
$ws.Range("E44").Value = "  -2.06%  "

$ws.Range("D45").Value = "'0.269"
$ws.Range("E45").Value = "  -7.85%  "

$ws.Range("D46").Value = "2.733.31"
$ws.Range("E46").Value = "  -4.03%  "

$ws.Range("D47").Value = "'37.44"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("D48").Value = "'129.18"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.19"
$ws.Range("E51").Value = "  -2.17%  "
